$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "28-02-2023 04:12"
$ws.Range("B10").Value = "hola"
